$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update inserts two new price rows (new reporting date 2023-04-05,
# serial 45021) right after the current row 516 ("Feria Lagunitas de Puerto
# Montt" / Pomelo data block), pushing the former rows 517:549 down to
# 519:551 (dimension grows from A1:T549 to A1:T551).
$ws.Range("A517:A518").EntireRow.Insert()

# New row 517 ("Primera" quality)
$ws.Cells.Item(517, 1).Value2 = 4
$ws.Cells.Item(517, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(517, 3).Value2 = "Los Lagos"
$ws.Cells.Item(517, 4).Value2 = 45021
$ws.Cells.Item(517, 5).Value2 = 10
$ws.Cells.Item(517, 6).Value2 = "Fruta"
$ws.Cells.Item(517, 7).Value2 = 100102
$ws.Cells.Item(517, 8).Value2 = "Cítricos"
$ws.Cells.Item(517, 9).Value2 = 100102006
$ws.Cells.Item(517, 10).Value2 = "Pomelo"
$ws.Cells.Item(517, 11).Value2 = "Start Ruby"
$ws.Cells.Item(517, 12).Value2 = "Primera"
$ws.Cells.Item(517, 13).Value2 = 80
$ws.Cells.Item(517, 14).Value2 = 14000
$ws.Cells.Item(517, 15).Value2 = 15000
$ws.Cells.Item(517, 16).Value2 = 14500
$ws.Cells.Item(517, 17).Value2 = "$/caja 14 kilos empedrada"
$ws.Cells.Item(517, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(517, 19).Value2 = 1036
$ws.Cells.Item(517, 20).Value2 = 14

# New row 518 ("Segunda" quality)
$ws.Cells.Item(518, 1).Value2 = 4
$ws.Cells.Item(518, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(518, 3).Value2 = "Los Lagos"
$ws.Cells.Item(518, 4).Value2 = 45021
$ws.Cells.Item(518, 5).Value2 = 10
$ws.Cells.Item(518, 6).Value2 = "Fruta"
$ws.Cells.Item(518, 7).Value2 = 100102
$ws.Cells.Item(518, 8).Value2 = "Cítricos"
$ws.Cells.Item(518, 9).Value2 = 100102006
$ws.Cells.Item(518, 10).Value2 = "Pomelo"
$ws.Cells.Item(518, 11).Value2 = "Start Ruby"
$ws.Cells.Item(518, 12).Value2 = "Segunda"
$ws.Cells.Item(518, 13).Value2 = 40
$ws.Cells.Item(518, 14).Value2 = 12000
$ws.Cells.Item(518, 15).Value2 = 12000
$ws.Cells.Item(518, 16).Value2 = 12000
$ws.Cells.Item(518, 17).Value2 = "$/caja 14 kilos empedrada"
$ws.Cells.Item(518, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(518, 19).Value2 = 857
$ws.Cells.Item(518, 20).Value2 = 14
